$wb = $excel.ActiveWorkbook

# --- Users sheet: remove "lucien" row and rotate admin's password hash ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B2").Value = "scrypt:32768:8:1`$gvT7jvcAopzF4dpg`$3840b05fde6c3060242de559b869dca1aeebf2b228cc9524caf70e8b78595b1cc6d67db69e34b3ab89f18d6226129fa33894f66919fce3caa5bfcff67cbf7b6f"
$wsUsers.Rows.Item(3).Delete()

# --- Directory sheet: supplier name correction ---
$wsDirectory = $wb.Worksheets.Item("Directory")
$wsDirectory.Range("C2").Value = "Lucien"

# --- Categories sheet: drop the "Spare" row ---
$wsCategories = $wb.Worksheets.Item("Categories")
$wsCategories.Rows.Item(2).Delete()

# --- Add the new "Vessels" sheet after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$wsVessels = $wb.Worksheets.Add($null, $wb.Worksheets.Item($sheetCount))
$wsVessels.Name = "Vessels"
$wsVessels.Range("A1").Value = "id"
$wsVessels.Range("B1").Value = "name"
$wsVessels.Range("A2").Value = 1
$wsVessels.Range("B2").Value = "Haykal"
